$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.604.33"
$ws.Range("E2").Value = "  -2.19%  "

# Row 3
$ws.Range("D3").Value = "1.842.71"
$ws.Range("E3").Value = "  -1.24%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.40%  "

# Row 5
$ws.Range("D5").Value = "'314.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("D7").Value = "'0.4243"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.36%  "

# Row 8
$ws.Range("D8").Value = "'0.3641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "

# Row 9
$ws.Range("D9").Value = "'45.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

# Row 10
$ws.Range("D10").Value = "'0.07265"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "

# Row 11
$ws.Range("D11").Value = "'0.8924"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.51%  "

# Row 12
$ws.Range("D12").Value = "'20.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.20%  "

# Row 13
$ws.Range("D13").Value = "1.817.45"
$ws.Range("E13").Value = "  -4.49%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.567"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.01%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.360"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.00%  "

# Row 16
$ws.Range("D16").Value = "'0.06872"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "'78.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.10%  "

# Row 19
$ws.Range("D19").Value = "'0.000008804"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.44%  "

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "

# Row 21
$ws.Range("E21").Value = "  -2.84%  "

# Row 22
$ws.Range("D22").Value = "27.595.76"
$ws.Range("E22").Value = "  -2.19%  "

# Row 23
$ws.Range("E23").Value = "  -2.29%  "

# Row 24
$ws.Range("E24").Value = "  -2.09%  "

# Row 25
$ws.Range("D25").Value = "2.049.13"
$ws.Range("E25").Value = "  -3.75%  "

# Row 26
$ws.Range("D26").Value = "'2.031"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

# Row 27
$ws.Range("D27").Value = "'155.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").Value = "'18.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "

# Row 29
$ws.Range("D29").Value = "'5.238"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30
$ws.Range("D30").Value = "'118.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.37%  "

# Row 31
$ws.Range("D31").Value = "'1.830"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.13%  "

# Row 32
$ws.Range("D32").Value = "'0.08894"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

# Row 33
$ws.Range("D33").Value = "'0.7770"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.21%  "

# Row 34
$ws.Range("D34").Value = "'4.571"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.52%  "

# Row 35
$ws.Range("D35").Value = "'2.958"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "

# Row 36
$ws.Range("D36").Value = "'1.104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.88%  "

# Row 37
$ws.Range("D37").Value = "'0.9998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "

# Row 38
$ws.Range("D38").Value = "'0.05412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("D39").Value = "'1.097"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.77%  "

# Row 40
$ws.Range("D40").Value = "'0.01920"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.49%  "

# Row 41
$ws.Range("D41").Value = "'2.769"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.59%  "

# Row 42
$ws.Range("D42").Value = "'6.836"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "

# Row 43
$ws.Range("D43").Value = "'0.5060"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.72%  "

# Row 44
$ws.Range("D44").Value = "'0.1652"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "

# Row 45
$ws.Range("D45").Value = "'8.199"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.89%  "

# Row 46
$ws.Range("D46").Value = "'0.06619"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.76%  "

# Row 47
$ws.Range("D47").Value = "'10.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.91%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4695"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.56%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'105.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.55%  "

# Row 50
$ws.Range("D50").Value = "'0.9997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

# Row 51
$ws.Range("D51").Value = "'1.625"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.75%  "
